# Remove the standalone italic "1 Samuel" paragraph that sits right
# after the "1SA" Heading2 paragraph (book-code / book-title pair at the
# top of the "09" resource section). The paragraph mark is removed along
# with the text so the following paragraph (the lone-space paragraph)
# shifts up in its place - matching the upstream diff exactly.

$d = $word.ActiveDocument

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd([char]13)

    if ($text -eq "1 Samuel" -and $p.Range.Italic -eq -1) {
        # Disambiguate from the unrelated "1 Samuel" Heading2 later in the
        # doc by confirming this one immediately follows the "1SA" para.
        if ($i -gt 1) {
            $prevText = $d.Paragraphs.Item($i - 1).Range.Text.TrimEnd([char]13)
            if ($prevText -eq "1SA") {
                $target = $p
            }
        }
    }
}

if ($target -ne $null) {
    $delRange = $d.Range($target.Range.Start, $target.Range.End)
    $null = $delRange.Delete()
}
